# Adds the new "pelada" results rows (181-202) to the Jogadores sheet,
# matching the upstream commit "Add files via upload".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, player name (column A), then the values for
# columns C..L (Vitorias, Empate, Derrotas, Gols, Partidas, Tarde de
# Vitoria, La barca, Craque do Dia, Gols Sofridos, Melhor Goleiro).
# Column B (Pontos) is intentionally left blank, matching the rest of
# the sheet.
$rowsData = @(
  @(181, "Coxinha", 4,0,4,0,1,0,0,0,0,0),
  @(182, "Leandrinho", 4,0,4,2,1,0,0,0,0,0),
  @(183, "Fabinho", 4,0,4,5,1,0,0,0,0,0),
  @(184, "Marcelão", 4,0,4,0,1,0,0,0,0,0),
  @(185, "Adriano", 4,0,4,1,1,0,0,0,0,0),
  @(186, "Corinthiano", 3,2,3,1,1,0,0,0,0,0),
  @(187, "Victor", 3,2,3,2,1,0,0,0,0,0),
  @(188, "Juscielio", 3,2,3,1,1,0,0,0,0,0),
  @(189, "Caio", 3,2,3,0,1,0,0,0,0,0),
  @(190, "Lucas", 3,2,3,2,1,0,0,0,0,0),
  @(191, "David", 4,2,2,1,1,1,0,0,0,0),
  @(192, "Fernando", 4,2,2,3,1,1,0,0,0,0),
  @(193, "Deiverson", 4,2,2,1,1,1,0,0,0,0),
  @(194, "Eder", 4,2,2,0,1,1,0,0,0,0),
  @(195, "Douglas", 4,2,2,2,1,1,0,0,0,0),
  @(196, "Digão", 2,0,4,1,1,0,1,0,0,0),
  @(197, "Eduardo", 2,0,4,4,1,0,1,0,0,0),
  @(198, "Guinha", 2,0,4,1,1,0,1,0,0,0),
  @(199, "Leah", 2,0,4,1,1,0,1,0,0,0),
  @(200, "João", 2,0,4,0,1,0,1,0,0,0),
  @(201, "Matheus", 6,2,7,0,1,0,1,0,14,0),
  @(202, "Chelin", 7,2,6,2,1,1,0,0,13,0)
)

foreach ($rowData in $rowsData) {
  $r = $rowData[0]
  $ws.Cells.Item($r, 1).Value = $rowData[1]
  # columns C (3) through L (12) map to indices 2..11 of $rowData
  for ($col = 3; $col -le 12; $col++) {
    $ws.Cells.Item($r, $col).Value = $rowData[$col - 1]
  }
}

# Match the new view state saved in the workbook: header row frozen,
# scrolled so row 185 is at the top, with I193 selected.
$win = $excel.ActiveWindow
[void]($win.FreezePanes = $false)
$ws.Range("A2").Select() | Out-Null
[void]($win.FreezePanes = $true)
$win.SetTopLeftVisibleCell("A185") | Out-Null
$ws.Range("I193").Select() | Out-Null
